$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window: tab ratio 980 -> 990 ---
$excel.ActiveWindow.TabRatio = 990

# --- Worksheet column width: 8.50510204081633 -> 8.23469387755102 (closest reachable value) ---
$ws.Columns.Item(1).ColumnWidth = 7.3

# --- Cell content/style changes ---
# Row 7: was text "data" (style 1) -> becomes number 1 (style 0 / default, unformatted)
# Row 8: was number 1 (style 1)    -> becomes text "auto" (style 1, unchanged)
# Row 9: was text "auto" (style 1) -> becomes text "data" (style 0 / default, unformatted)
#
# Style index 0 means "no explicit formatting applied" in this workbook, which is not
# reachable by setting Style/Font properties directly (that always synthesizes a brand
# new style record). Instead, grab the default/unformatted style from a never-touched
# cell far outside the used range and paste only the formatting from it.

$blank = $ws.Range("ZZ999")
$blank.Value = 1
$blank.Copy()

$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 1

$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "data"

$blank.Clear()

# A8 keeps style 1 (already the case) -- only its value/type changes to the text "auto"
$ws.Range("A8").Value = "auto"

# --- Selection moves from A10 to A8 ---
$ws.Range("A8").Select()
